# Update column G ("K") values on Sheet1 per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 3
    3  = 0
    5  = 1
    6  = 3
    7  = 0
    8  = 3
    9  = 3
    10 = 2
    12 = 1
    13 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
